$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.088.78'
$ws.Range('E2').Value = '  +0.48%  '

$ws.Range('D3').Value = '2.759.18'
$ws.Range('E3').Value = '  +1.32%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.68%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.56'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.79%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.11%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.608'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.10%  '

$ws.Range('E9').Value = '  -1.47%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.70'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -14.65%  '

$ws.Range('E11').Value = '  -0.73%  '

$ws.Range('E12').Value = '  -2.01%  '

$ws.Range('D13').Value = '3.247.63'
$ws.Range('E13').Value = '  +1.25%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.91'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.36%  '

$ws.Range('D15').Value = '63.780.74'
$ws.Range('E15').Value = '  -0.02%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000152'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.12%  '

$ws.Range('D17').Value = '2.761.43'
$ws.Range('E17').Value = '  +0.77%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.17'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.96%  '

$ws.Range('E19').Value = '  +0.62%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '360.48'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.03%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.01%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.549'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.85%  '

$ws.Range('E23').Value = '  +0.38%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.81'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.31%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.171'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.49%  '

$ws.Range('E26').Value = '  +0.47%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.12%  '

$ws.Range('D28').Value = '0.0₃0927'
$ws.Range('E28').Value = '  +2.25%  '

$ws.Range('E29').Value = '  -1.30%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.03'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.07%  '

$ws.Range('E31').Value = '  +1.17%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '167.59'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.31%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.35'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.56%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.94'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.61%  '

$ws.Range('E35').Value = '  +0.11%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.46'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.10%  '

$ws.Range('E37').Value = '  +0.26%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.992'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.37%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.29'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +12.69%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.17'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.01%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '330.02'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.20%  '

$ws.Range('E42').Value = '  +0.36%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.63'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.22%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0596'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.96%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.77'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.21%  '

$ws.Range('E46').Value = '  +1.40%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.635'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.16%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '136.41'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.89%  '

$ws.Range('E49').Value = '  +0.91%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.27%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.05'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.65%  '
